# Fixed #348 Hyperlinks from sub-template does not work.
#
# The second paragraph of the document starts with a stray, completely
# empty run:
#
#   <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t/></w:r>
#
# We replace that empty run with a "_GoBack" bookmark pair
# (<w:bookmarkStart .../><w:bookmarkEnd .../>) placed at the exact same
# spot, immediately after the paragraph's <w:pPr> and before the field
# that follows it.

$d = $word.ActiveDocument

# Locate the target paragraph (2nd paragraph -> the one beginning with the
# empty run, followed by the " m: a + a + b " field).
$p = $d.Paragraphs.Item(2)
$start = $p.Range.Start

# The leading run is empty (zero-width), so a collapsed Range pointing at
# it can't be targeted for deletion directly. Give it real, addressable
# content first by writing into the collapsed range: Word reuses the
# existing (empty) run for this rather than minting a new one.
$markerRange = $d.Range($start, $start)
$markerRange.Text = "x"

# Now that the run has one real character, select it with a proper
# (non-collapsed) range and clear it. Because the range actually spans
# text, this cleanly removes the character together with the now-empty
# run, instead of behaving like a delete of the whole paragraph.
$markerRange2 = $d.Range($start, $start + 1)
$markerRange2.Text = ""

# Insert the _GoBack bookmark exactly where that empty run used to sit.
$bookmarkRange = $d.Range($start, $start)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
